# "9th Stab - Cosmetic Changes"
# Insert two new weekly-snapshot columns ("Jun_17", "Jun_15") to the left of
# the existing "Jun_13" / "Jun_10" columns, pushing the old data right by two
# columns, and backfill the new columns with the same "UN" placeholder used
# throughout the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at B:C - this shifts the former B (Jun_13) and
# C (Jun_10) columns to D and E respectively, carrying their values/styles
# with them.
$ws.Columns("B:C").Insert()

# Match the column width of the neighbouring data columns (~8 characters).
$ws.Columns("C").ColumnWidth = 7.14
$ws.Columns("D").ColumnWidth = 7.14

# New header row.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# New data columns: fill with the same "UN" placeholder used elsewhere.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}
